$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $ref, $val) {
    $rng = $ws.Range($ref)
    if ($val -match '^[+-]?[0-9]*\.?[0-9]+([eE][+-]?[0-9]+)?$') {
        $s = $rng.Style
        $rng.NumberFormat = "@"
        $rng.Value = $val
        $rng.Style = $s
    } else {
        $rng.Value = $val
    }
}

Set-TextValue $ws "D2" "65.268.40"
Set-TextValue $ws "E2" "  -0.70%  "
Set-TextValue $ws "D3" "3.541.01"
Set-TextValue $ws "E3" "  +2.71%  "
Set-TextValue $ws "D4" "1.00"
Set-TextValue $ws "E4" "  +0.05%  "
Set-TextValue $ws "D5" "603.44"
Set-TextValue $ws "E5" "  +1.69%  "
Set-TextValue $ws "D6" "139.66"
Set-TextValue $ws "E6" "  +1.77%  "
Set-TextValue $ws "D7" "3.540.77"
Set-TextValue $ws "E7" "  +2.74%  "
Set-TextValue $ws "E8" "  +0.13%  "
Set-TextValue $ws "D9" "0.491"
Set-TextValue $ws "E9" "  -1.92%  "
Set-TextValue $ws "D10" "0.125"
Set-TextValue $ws "E10" "  +2.27%  "
Set-TextValue $ws "D11" "6.96"
Set-TextValue $ws "E11" "  -5.98%  "
Set-TextValue $ws "D12" "0.392"
Set-TextValue $ws "E12" "  +3.08%  "
Set-TextValue $ws "D13" "4.147.30"
Set-TextValue $ws "E13" "  +3.13%  "
Set-TextValue $ws "D14" "0.0000187"
Set-TextValue $ws "E14" "  +2.15%  "
Set-TextValue $ws "D15" "27.33"
Set-TextValue $ws "E15" "  +3.04%  "
Set-TextValue $ws "D16" "3.551.73"
Set-TextValue $ws "E16" "  +2.83%  "
Set-TextValue $ws "E17" "  +1.66%  "
Set-TextValue $ws "D18" "65.382.20"
Set-TextValue $ws "E18" "  -0.35%  "
Set-TextValue $ws "D19" "10.36"
Set-TextValue $ws "E19" "  +4.64%  "
Set-TextValue $ws "D20" "5.96"
Set-TextValue $ws "E20" "  +1.46%  "
Set-TextValue $ws "D21" "14.38"
Set-TextValue $ws "E21" "  +4.81%  "
Set-TextValue $ws "D22" "395.92"
Set-TextValue $ws "E22" "  +0.58%  "
Set-TextValue $ws "D23" "0.576"
Set-TextValue $ws "E23" "  +3.88%  "
Set-TextValue $ws "D24" "3.684.01"
Set-TextValue $ws "E24" "  +2.81%  "
Set-TextValue $ws "D25" "73.99"
Set-TextValue $ws "E25" "  +0.74%  "
Set-TextValue $ws "D26" "0.999"
Set-TextValue $ws "E26" "  -0.18%  "
Set-TextValue $ws "E27" "  +9.23%  "
Set-TextValue $ws "D28" "7.83"
Set-TextValue $ws "E28" "  +8.37%  "
Set-TextValue $ws "E29" "  +0.05%  "
Set-TextValue $ws "E30" "  +2.13%  "
Set-TextValue $ws "D31" "8.31"
Set-TextValue $ws "E31" "  +0.51%  "
Set-TextValue $ws "D32" "3.554.98"
Set-TextValue $ws "E32" "  +3.08%  "
Set-TextValue $ws "E33" "  +0.00%  "
Set-TextValue $ws "D34" "23.82"
Set-TextValue $ws "E34" "  +3.40%  "
Set-TextValue $ws "D35" "0.146"
Set-TextValue $ws "E35" "  -0.41%  "
Set-TextValue $ws "D36" "1.30"
Set-TextValue $ws "E36" "  +9.70%  "
Set-TextValue $ws "D37" "7.00"
Set-TextValue $ws "E37" "  +0.43%  "
Set-TextValue $ws "D38" "1.56"
Set-TextValue $ws "E38" "  +3.84%  "
Set-TextValue $ws "D39" "168.96"
Set-TextValue $ws "E39" "  -1.97%  "
Set-TextValue $ws "D40" "4.98"
Set-TextValue $ws "E40" "  +3.17%  "
Set-TextValue $ws "D41" "0.0815"
Set-TextValue $ws "E41" "  +6.13%  "
Set-TextValue $ws "D42" "0.829"
Set-TextValue $ws "E42" "  +0.22%  "
Set-TextValue $ws "D43" "26.79"
Set-TextValue $ws "E43" "  +16.18%  "
Set-TextValue $ws "D44" "43.01"
Set-TextValue $ws "E44" "  -1.59%  "
Set-TextValue $ws "E45" "  +0.12%  "
Set-TextValue $ws "D46" "4.45"
Set-TextValue $ws "E46" "  +0.50%  "
Set-TextValue $ws "E47" "  +9.58%  "
Set-TextValue $ws "D48" "1.68"
Set-TextValue $ws "E48" "  +3.69%  "
Set-TextValue $ws "D49" "2.457.76"
Set-TextValue $ws "E49" "  +11.20%  "
Set-TextValue $ws "E50" "  +3.41%  "
Set-TextValue $ws "B51" "Bittensor"
Set-TextValue $ws "C51" "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
Set-TextValue $ws "D51" "304.17"
Set-TextValue $ws "E51" "  +8.77%  "
